$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: id=5, nome=GABRIEL BONARETTI
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "GABRIEL BONARETTI"
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0

# Row 8: id=6, nome=PEDRO
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "PEDRO"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
